$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.100.64"
$ws.Range("E2").Value = "  -4.33%  "

# Row 3
$ws.Range("D3").Value = "1.835.76"
$ws.Range("E3").Value = "  -2.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.04%  "

# Row 6
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.78%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3872"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.54%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07894"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.44%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9607"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.25%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.659"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.30%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.896"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.74%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.758.41"
$ws.Range("E15").Value = "  -9.85%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06815"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.64%  "

# Row 17
$ws.Range("E17").Value = "  -0.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009990"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.73%  "

# Row 21
$ws.Range("E21").Value = "  -0.20%  "

# Row 22
$ws.Range("D22").Value = "28.131.71"
$ws.Range("E22").Value = "  -4.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.328"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.78%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.100"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.14%  "

# Row 26
$ws.Range("D26").Value = "2.132.02"
$ws.Range("E26").Value = "  -1.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.57%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.693"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.979"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.90%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09295"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.10%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9308"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.303"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.03%  "

# Row 35
$ws.Range("E35").Value = "  -4.54%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.347"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.93%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05892"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02152"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.145"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.725"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5596"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.70%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.926"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1769"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.229"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.06%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.82%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.178"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5270"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07011"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.837"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.94%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.01%  "

# Row 51
$ws.Range("E51").Value = "  -0.10%  "
